# Insert a new data row after the header/existing rows at row 76, shifting
# all subsequent rows (old 76..196) down to (77..197), then populate the
# newly inserted row 76 with a new price-list entry for Membrillo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 76 (pushes old row 76 -> 77, etc.)
$ws.Rows.Item(76).Insert()

# Copy the (now shifted-down) row 77 static attributes into the new row 76,
# since the new record shares the same market/product/quality metadata.
for ($c = 1; $c -le 20; $c++) {
    $ws.Cells.Item(76, $c).Value2 = $ws.Cells.Item(77, $c).Value2
}

# Overwrite the columns that differ for this new entry:
#   D = Fecha (date serial), N/O/P = precios, S = precio $/Kg
$ws.Cells.Item(76, 4).Value2 = 44757
$ws.Cells.Item(76, 14).Value2 = 10000
$ws.Cells.Item(76, 15).Value2 = 10000
$ws.Cells.Item(76, 16).Value2 = 10000
$ws.Cells.Item(76, 19).Value2 = 556
